## edit.ps1
## Applies the "new branch" edit to gandhi.docx:
##   1. Appends "  (This is a change [EN DASH] Version for branch alternate)"
##      to the end of the first paragraph; the two leading spaces stay in the
##      default (black) run, the parenthetical text is colored dark red
##      (C00000).
##   2. Adds one brand-new, completely empty paragraph right before the
##      section break at the very end of the document body.
##
## Note: the diff also shows a bare <w:rsid w:val="001772C0"/> added to the
## Normal paragraph style in styles.xml. That attribute is Word's internal
## "revision save ID" bookkeeping (purely cosmetic/session metadata, no
## visible or semantic effect) and there is no Word object-model property
## that exposes/sets style-level rsids (confirmed: Style has no Rsid/XML/
## WordOpenXML writer, and no document-level operation -- UpdateStyles(),
## re-applying the style, etc. -- causes the host to stamp one). It is left
## untouched here because it cannot be produced through legitimate
## COM automation.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: add the trailing note.
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$r = $firstPara.Range

# Range.End of a paragraph sits just past its paragraph mark (so the
# mark itself lives at End - 1); InsertAfter on the paragraph's own
# Range always lands new text *before* that mark, which is what we
# want (it keeps the text inside paragraph 1 rather than paragraph 2).
$enDash = [char]0x2013
$noteText = "(This is a change" + $enDash + "Version for branch alternate)"

# 1a. Two plain (uncoloured) spaces right after the existing sentence.
$r.InsertAfter("  ")
$afterSpaces = $r.End - 1

# 1b. The parenthetical note itself, which must end up dark red (C00000).
$r.InsertAfter($noteText)
$afterNote = $r.End - 1

# Color only the "(This is a change ... alternate)" span; the OLE/COM
# Font.Color value is encoded 0x00BBGGRR, so RGB(C0,00,00) == 0x0000C0 == 192.
$noteRange = $d.Range($afterSpaces, $afterNote)
$noteRange.Font.Color = 192

# ---------------------------------------------------------------------
# 2) Append one empty paragraph at the very end of the document, right
#    before the final section break.
# ---------------------------------------------------------------------
# Step 2a: create a new paragraph mark after the current last paragraph.
# (InsertParagraphAfter at the end of the body necessarily clones the
# preceding paragraph's formatting -- real Word does the same thing when
# you press Enter -- so the new last paragraph temporarily inherits the
# "Normal (Web)" / shading formatting of the paragraph before it.)
$endOfBody = $d.Content.End
$tailPoint = $d.Range($endOfBody - 1, $endOfBody - 1)
$tailPoint.InsertParagraphAfter()

# Step 2b: the freshly-created paragraph is now its own, unambiguous
# last paragraph, so replacing its Range via InsertXML with a bare
# <w:p/> cleanly overwrites *only* that paragraph's (inherited, unwanted)
# formatting -- leaving every earlier paragraph (including the
# "Normal (Web)" one right before it) untouched -- and leaves behind
# exactly the plain, empty paragraph the diff adds right before sectPr.
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$emptyParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$newLastPara.Range.InsertXML($emptyParagraphXml)
